$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark attendance "O" for these cells (same as other subjects already filled in column G)
$ws.Range("G2").Value = "O"
$ws.Range("G3").Value = "O"
$ws.Range("G5").Value = "O"
$ws.Range("G6").Value = "O"
$ws.Range("G7").Value = "O"
$ws.Range("G9").Value = "O"

# Move the active selection to G10
$ws.Range("G10").Select()
